$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -13.625
$ws.Range("C4").Value = -12.041
$ws.Range("B7").Value = 5.254
$ws.Range("A8").Value = -22.151
$ws.Range("A10").Value = -21.431
$ws.Range("C11").Value = -12.908
$ws.Range("A12").Value = -21.768
$ws.Range("B14").Value = 5.208
$ws.Range("C14").Value = -11.869
$ws.Range("B15").Value = 5.072000000000001
$ws.Range("A18").Value = -21.514
$ws.Range("B18").Value = 6.989
$ws.Range("C18").Value = -11.657
$ws.Range("C19").Value = -11.871
$ws.Range("B20").Value = 6.489999999999999
$ws.Range("C21").Value = -12.016
$ws.Range("A25").Value = -21.761
$ws.Range("C27").Value = -13.469
$ws.Range("B29").Value = 5.699
$ws.Range("B30").Value = 5.552
$ws.Range("B31").Value = 5.396
$ws.Range("C31").Value = -12.903
$ws.Range("B35").Value = 7.45
$ws.Range("A37").Value = -20.727
$ws.Range("C38").Value = -13.055
$ws.Range("B40").Value = 8.597999999999999
$ws.Range("C42").Value = -12.42
$ws.Range("B44").Value = 5.954000000000001
$ws.Range("C44").Value = -13.321
$ws.Range("C47").Value = -12.616
$ws.Range("B50").Value = 5.24
$ws.Range("B54").Value = 5.067
$ws.Range("A55").Value = -21.838
$ws.Range("C56").Value = -13.643
$ws.Range("C58").Value = -13.27
$ws.Range("C65").Value = -12.282
$ws.Range("A68").Value = -21.854
$ws.Range("B68").Value = 5.422
$ws.Range("C73").Value = -12.255
$ws.Range("B76").Value = 5.971
$ws.Range("A77").Value = -20.729
$ws.Range("A78").Value = -20.562
$ws.Range("A79").Value = -21.366
$ws.Range("A80").Value = -20.442
$ws.Range("A81").Value = -22.001
$ws.Range("A82").Value = -22.016
$ws.Range("A84").Value = -21.682
$ws.Range("B87").Value = 4.853000000000001
$ws.Range("B88").Value = 5.16
$ws.Range("C90").Value = -13.331
$ws.Range("B92").Value = 6.272
$ws.Range("C92").Value = -10.91
$ws.Range("C94").Value = -10.505
$ws.Range("C95").Value = -11.451
$ws.Range("B96").Value = 6.531000000000001
$ws.Range("B98").Value = 5.615
$ws.Range("A101").Value = -21.006
$ws.Range("B101").Value = 6.622
$ws.Range("C101").Value = -12.136
$ws.Range("A102").Value = -21.32
$ws.Range("B102").Value = 6.409999999999999
